# Correct the Copa America 2024 schedule/matchups on the "Matches" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matches")

# --- Match 2 (row 3): fix kickoff time only ---
$ws.Range("D3").Value = 45465

# --- Match 3 (row 4): swap teams, fix time, renumber venue ---
$ws.Range("B4").Value = "B1"
$ws.Range("C4").Value = "B4"
$ws.Range("D4").Value = 45466.041666666672
$ws.Range("E4").Value = 4

# --- Match 4 (row 5): swap teams, fix time, renumber venue ---
$ws.Range("B5").Value = "B2"
$ws.Range("C5").Value = "B3"
$ws.Range("D5").Value = 45465.916666666664
$ws.Range("E5").Value = 3

# --- Match 5 (row 6): fix kickoff time only ---
$ws.Range("D6").Value = 45466.916666666672

# --- Match 6 (row 7): fix kickoff time only ---
$ws.Range("D7").Value = 45467.041666666664

# --- Match 7 (row 8): swap teams, fix time, renumber venue ---
$ws.Range("B8").Value = "D1"
$ws.Range("C8").Value = "D4"
$ws.Range("D8").Value = 45468.041666666664
$ws.Range("E8").Value = 6

# --- Match 8 (row 9): swap teams, fix time, renumber venue ---
$ws.Range("B9").Value = "D2"
$ws.Range("C9").Value = "D3"
$ws.Range("D9").Value = 45467.916666666672
$ws.Range("E9").Value = 4

# --- Match 9 (row 10): swap teams, fix time, renumber venue ---
$ws.Range("B10").Value = "A3"
$ws.Range("C10").Value = "A1"
$ws.Range("D10").Value = 45469.041666666664
$ws.Range("E10").Value = 8

# --- Match 10 (row 11): swap teams, fix time, renumber venue ---
$ws.Range("B11").Value = "A2"
$ws.Range("C11").Value = "A4"
$ws.Range("D11").Value = 45468.916666666672
$ws.Range("E11").Value = 7

# --- Match 11 (row 12): swap teams, renumber venue (time unchanged) ---
$ws.Range("B12").Value = "B3"
$ws.Range("C12").Value = "B1"
$ws.Range("E12").Value = 6

# --- Match 12 (row 13): swap teams, fix time, renumber venue ---
$ws.Range("B13").Value = "B2"
$ws.Range("C13").Value = "B4"
$ws.Range("D13").Value = 45469.916666666664
$ws.Range("E13").Value = 9

# --- Match 14 (row 15): fix kickoff time only ---
$ws.Range("D15").Value = 45471.041666666664

# --- Match 15 (row 16): swap teams, renumber venue (time unchanged) ---
$ws.Range("B16").Value = "D3"
$ws.Range("C16").Value = "D1"
$ws.Range("E16").Value = 9

# --- Match 16 (row 17): swap teams, fix time, renumber venue ---
$ws.Range("B17").Value = "D2"
$ws.Range("C17").Value = "D4"
$ws.Range("D17").Value = 45471.916666666664
$ws.Range("E17").Value = 10

# --- Match 19 (row 20): swap teams, fix time, renumber venue ---
$ws.Range("B20").Value = "B1"
$ws.Range("C20").Value = "B2"
$ws.Range("D20").Value = 45474
$ws.Range("E20").Value = 10

# --- Match 20 (row 21): swap teams, fix time, renumber venue ---
$ws.Range("B21").Value = "B4"
$ws.Range("C21").Value = "B3"
$ws.Range("D21").Value = 45474
$ws.Range("E21").Value = 12

# --- Match 21 (row 22): swap teams, fix time, renumber venue ---
$ws.Range("B22").Value = "C1"
$ws.Range("C22").Value = "C2"
$ws.Range("D22").Value = 45475.041666666672
$ws.Range("E22").Value = 13

# --- Match 22 (row 23): swap teams, fix time, renumber venue ---
$ws.Range("B23").Value = "C4"
$ws.Range("C23").Value = "C3"
$ws.Range("D23").Value = 45475.041666666664
$ws.Range("E23").Value = 11

# --- Match 23 (row 24): swap teams, fix time, renumber venue ---
$ws.Range("B24").Value = "D1"
$ws.Range("C24").Value = "D2"
$ws.Range("D24").Value = 45476.041666666664
$ws.Range("E24").Value = 3

# --- Match 24 (row 25): swap teams, fix time, renumber venue ---
$ws.Range("B25").Value = "D4"
$ws.Range("C25").Value = "D3"
$ws.Range("D25").Value = 45476.041666666672
$ws.Range("E25").Value = 12

# --- Match 30 (row 31): fix date ---
$ws.Range("D31").Value = 45484

# Re-apply the table's sort state (sorted by "match" column) to match the
# corrected row order used when the schedule was fixed.
$tbl = $ws.ListObjects.Item("matches")
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add2($ws.Range("A2:A33"), $null, 1, $null, 0)
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# Restore the selected cell shown when the workbook was last saved.
$ws.Range("E32").Select()
